# ---------------------------------------------------------------------------
# feat(formulaire): updated ligne-formulaire update to delete links to
# formulaire.
#
# - extend the "area"/"each" jxls comments so their lastCell points at J5
#   instead of J4 (the template grew one row)
# - add a new "jx:each(... section.lignesFormulaire ...)" comment on A5
# - add a new template row (row 5) with ${ligne.code} / ${ligne.libelle} /
#   ${ligne.contenu} placeholders, merged C5:F5 and G5:I5
# - the old A1:J1 title merge is removed -- the title row is now made up of
#   individually bordered/filled cells instead of a single merged cell
# - every cell switches from centered (or default) alignment to left
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- xlsx constants used below -------------------------------------------
$xlEdgeLeft     = 7
$xlEdgeTop      = 8
$xlEdgeBottom   = 9
$xlEdgeRight    = 10
$xlContinuous   = 1
$xlLineStyleNone = -4142
$xlLeft         = -4131

# ---------------------------------------------------------------------------
# 1) Comments: bump lastCell="J4" -> lastCell="J5" on the four existing
#    comments, and add the new A5 comment for the "ligne" loop.
# ---------------------------------------------------------------------------
$ws.Range("A1").Comment.Text("Auteur:`njx:area(lastCell=`"J5`")")
$ws.Range("A2").Comment.Text("Auteur:`njx:each(items=`"formulaires`", var=`"formulaire`", multisheet=`"sheetNames`", lastCell=`"J5`")")
$ws.Range("A3").Comment.Text("Auteur:`njx:each(items=`"formulaire.sections`", var=`"section`", lastCell=`"J5`")")
$ws.Range("A4").Comment.Text("Auteur:`njx:each(items=`"section.sections`" var=`"section`", lastCell=`"J5`")")

$newComment = $ws.Range("A5").AddComment()
$newComment.Text("Auteur:`njx:each(items=`"section.lignesFormulaire`" var=`"ligne`" lastCell=`"J5`")")

# ---------------------------------------------------------------------------
# 2) Title row: drop the A1:J1 merge, keep the same fill but rebuild the
#    border as a run of individually-bordered cells (left cell keeps the
#    left edge, right cell keeps the right edge, middle cells only keep
#    top/bottom).
# ---------------------------------------------------------------------------
$ws.Range("A1:J1").UnMerge()

$titleCols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $titleCols) {
  $cell = $ws.Range($col + "1")
  $cell.Borders($xlEdgeTop).LineStyle = $xlContinuous
  $cell.Borders($xlEdgeBottom).LineStyle = $xlContinuous
  if ($col -eq "A") {
    $cell.Borders($xlEdgeLeft).LineStyle = $xlContinuous
    $cell.Borders($xlEdgeRight).LineStyle = $xlLineStyleNone
  } elseif ($col -eq "J") {
    $cell.Borders($xlEdgeLeft).LineStyle = $xlLineStyleNone
    $cell.Borders($xlEdgeRight).LineStyle = $xlContinuous
  } else {
    $cell.Borders($xlEdgeLeft).LineStyle = $xlLineStyleNone
    $cell.Borders($xlEdgeRight).LineStyle = $xlLineStyleNone
  }
}

# ---------------------------------------------------------------------------
# 2b) "${formulaire.name}" merged box (B2:F2): same left/right/top/bottom ->
#     segmented top+bottom-everywhere border treatment as the title row
#     (previously this merge only had left/right/bottom, no top).
# ---------------------------------------------------------------------------
$nameCols = @("B","C","D","E","F")
foreach ($col in $nameCols) {
  $cell = $ws.Range($col + "2")
  $cell.Borders($xlEdgeTop).LineStyle = $xlContinuous
  $cell.Borders($xlEdgeBottom).LineStyle = $xlContinuous
  if ($col -eq "B") {
    $cell.Borders($xlEdgeLeft).LineStyle = $xlContinuous
    $cell.Borders($xlEdgeRight).LineStyle = $xlLineStyleNone
  } elseif ($col -eq "F") {
    $cell.Borders($xlEdgeLeft).LineStyle = $xlLineStyleNone
    $cell.Borders($xlEdgeRight).LineStyle = $xlContinuous
  } else {
    $cell.Borders($xlEdgeLeft).LineStyle = $xlLineStyleNone
    $cell.Borders($xlEdgeRight).LineStyle = $xlLineStyleNone
  }
}

# ---------------------------------------------------------------------------
# 3) New template row 5: values + merges.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = '${ligne.code}'
$ws.Range("C5").Value = '${ligne.libelle}'
$ws.Range("G5").Value = '${ligne.contenu}'

$ws.Range("C5:F5").Merge()
$ws.Range("G5:I5").Merge()

# ---------------------------------------------------------------------------
# 4) Alignment: the whole template switches from center/default to left.
#    Touch exactly the cells that already carry content/formatting so we
#    don't materialize stray empty cells outside the template's footprint.
# ---------------------------------------------------------------------------
$ws.Range("A1:J2").HorizontalAlignment = $xlLeft
$ws.Range("A3:I3").HorizontalAlignment = $xlLeft
$ws.Range("B4:I4").HorizontalAlignment = $xlLeft
$ws.Range("B5:I5").HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------------
# 5) Cosmetic: move the active selection like the authored file does.
# ---------------------------------------------------------------------------
$ws.Range("F10").Select()
